$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value to a cell without Excel coercing
# numeric-looking strings (e.g. "212.20", "0.536") into floating point
# numbers -- forces text format for the write, then restores the cells
# original style so no spurious style index is introduced.
function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

$ws.Range("D2").Value = "27.429.30"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "1.643.57"
$ws.Range("E3").Value = "  -1.36%  "
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.06%  "
Set-TextValue $ws.Range("D5") "212.20"
$ws.Range("E5").Value = "  -1.45%  "
Set-TextValue $ws.Range("D6") "0.536"
$ws.Range("E6").Value = "  +4.25%  "
$ws.Range("E7").Value = "  -0.05%  "
Set-TextValue $ws.Range("D8") "23.25"
$ws.Range("E8").Value = "  -1.47%  "
Set-TextValue $ws.Range("D9") "0.257"
$ws.Range("E9").Value = "  -2.37%  "
$ws.Range("E10").Value = "  -2.01%  "
Set-TextValue $ws.Range("D11") "0.0890"
$ws.Range("E11").Value = "  +0.72%  "
$ws.Range("D12").Value = "1.874.84"
$ws.Range("E12").Value = "  -1.42%  "
$ws.Range("D13").Value = "1.628.85"
$ws.Range("E13").Value = "  -3.76%  "
Set-TextValue $ws.Range("D14") "4.04"
$ws.Range("E14").Value = "  -2.84%  "
Set-TextValue $ws.Range("D15") "0.559"
$ws.Range("E15").Value = "  +0.29%  "
Set-TextValue $ws.Range("D16") "64.29"
$ws.Range("E16").Value = "  -2.93%  "
$ws.Range("D17").Value = "27.400.98"
$ws.Range("E17").Value = "  -0.72%  "
Set-TextValue $ws.Range("D18") "228.80"
$ws.Range("E18").Value = "  -8.28%  "
$ws.Range("E19").Value = "  -1.79%  "
Set-TextValue $ws.Range("D20") "7.52"
$ws.Range("E20").Value = "  -0.37%  "
$ws.Range("E21").Value = "  -0.02%  "
Set-TextValue $ws.Range("D22") "4.32"
$ws.Range("E22").Value = "  -4.07%  "
Set-TextValue $ws.Range("D23") "9.34"
$ws.Range("E23").Value = "  +0.35%  "
Set-TextValue $ws.Range("D24") "2.04"
$ws.Range("E24").Value = "  +0.02%  "
Set-TextValue $ws.Range("D25") "148.01"
$ws.Range("E25").Value = "  +1.32%  "
Set-TextValue $ws.Range("D26") "0.114"
$ws.Range("E26").Value = "  +2.42%  "
Set-TextValue $ws.Range("D27") "6.94"
$ws.Range("E27").Value = "  -2.94%  "
$ws.Range("E28").Value = "  -0.05%  "
Set-TextValue $ws.Range("D29") "15.53"
$ws.Range("E29").Value = "  -5.70%  "
$ws.Range("E30").Value = "  -4.54%  "
Set-TextValue $ws.Range("D31") "0.0488"
$ws.Range("E31").Value = "  -3.95%  "
$ws.Range("E32").Value = "  -2.54%  "
Set-TextValue $ws.Range("D33") "3.10"
$ws.Range("E33").Value = "  -0.36%  "
$ws.Range("D34").Value = "1.411.05"
$ws.Range("E34").Value = "  -4.24%  "
$ws.Range("E35").Value = "  -0.87%  "
$ws.Range("E36").Value = "  -0.26%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D37") "0.880"
$ws.Range("E37").Value = "  -6.25%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D38") "0.561"
$ws.Range("E38").Value = "  -2.38%  "
Set-TextValue $ws.Range("D39") "0.0166"
$ws.Range("E39").Value = "  -3.43%  "
$ws.Range("E40").Value = "  +0.23%  "
$ws.Range("E41").Value = "  -0.05%  "
Set-TextValue $ws.Range("D42") "2.49"
$ws.Range("E42").Value = "  -0.91%  "
Set-TextValue $ws.Range("D43") "5.48"
$ws.Range("E43").Value = "  +1.13%  "
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D45") "64.60"
$ws.Range("E45").Value = "  -7.21%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D46") "0.789"
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("D47").Value = "1.785.68"
$ws.Range("E47").Value = "  -1.30%  "
Set-TextValue $ws.Range("D48") "1.65"
$ws.Range("E48").Value = "  -3.51%  "
Set-TextValue $ws.Range("D49") "87.45"
$ws.Range("E49").Value = "  -2.20%  "
$ws.Range("E50").Value = "  -2.93%  "
Set-TextValue $ws.Range("D51") "0.0987"
$ws.Range("E51").Value = "  -3.34%  "
